## Sprint2Backlog.xlsx - "burndown and retrospective report added"
##
## Real content changes applied:
##  1. Row 23 ("Sprint 2" sheet) is turned from a plain backlog-item row into a
##     new section header row (matching the style used by rows 12, 15, 19, 26,
##     29, 32, 35, 38, 41, 44):
##       - A23 text updated to the new story title
##       - A23 / B23 formatting swapped to the "section header" look
##       - F23 gains the section's day-5 total (3)
##  2. Row 47 ("Total" row) gains the burndown totals across F:L (one literal
##     count plus SUM() formulas mirroring the existing G:L pattern).
##  3. Selection/viewport updated to match where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")

# --- 1. Row 23 becomes the new "16. ... service provider information page" header ---

# Clone the look of an existing section-header row (row 12: grey fill + full
# border) onto the A:B and F cells of row 23 before touching the values.
$ws.Range("A12:B12").Copy()
$ws.Range("A23:B23").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F12").Copy()
$ws.Range("F23").PasteSpecial(-4122)       # xlPasteFormats

$ws.Range("A23").Value = "16. As a customer, I can view service provider information page"
$ws.Range("B23").Value = 4
$ws.Range("F23").Value = 3

# --- 2. Row 47 "Total" row gains the burndown totals ---

$ws.Range("F47").Value = 48
$ws.Range("G47").Formula = "=SUM(G3:G25)"
$ws.Range("H47:L47").Formula = "=SUM(H3:H25)"

# --- 3. Restore the viewport/selection left behind by the author ---

$ws.Range("G48").Select()
